# Apply scheduled-runner updates to Sheets (commit: chore: update Sheets via scheduled runner)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2165805.5
$ws.Range("J17").Value = 2274028.8
$ws.Range("L17").Value = 6822086.399999999
$ws.Range("N17").Value = -6822422.399999999
$ws.Range("H28").Value = 321.08334
$ws.Range("I28").Value = 338.8
$ws.Range("K28").Value = 338.8
$ws.Range("M28").Value = 146.2
$ws.Range("H33").Value = 525.375
$ws.Range("I33").Value = 286.54544
$ws.Range("K33").Value = 286.54544
$ws.Range("M33").Value = -57.54543999999999
$ws.Range("H40").Value = 4696.857
$ws.Range("I40").Value = 2495
$ws.Range("J40").Value = 5577.6
$ws.Range("K40").Value = 2495
$ws.Range("L40").Value = 5577.6
$ws.Range("M40").Value = -2320
$ws.Range("N40").Value = -5927.6
$ws.Range("H69").Value = 9798.154
$ws.Range("I69").Value = 7274.625
$ws.Range("J69").Value = 13835.8
$ws.Range("K69").Value = 21823.875
$ws.Range("L69").Value = 41507.39999999999
$ws.Range("M69").Value = -20949.875
$ws.Range("N69").Value = -43255.39999999999
$ws.Range("H72").Value = 9798.154
$ws.Range("I72").Value = 7274.625
$ws.Range("J72").Value = 13835.8
$ws.Range("K72").Value = 65471.625
$ws.Range("L72").Value = 124522.2
$ws.Range("M72").Value = -61103.625
$ws.Range("N72").Value = -133258.2
$ws.Range("H88").Value = 1395.1428
$ws.Range("I88").Value = 1115
$ws.Range("J88").Value = 1550.7778
$ws.Range("K88").Value = 1115
$ws.Range("L88").Value = 1550.7778
$ws.Range("M88").Value = -709
$ws.Range("N88").Value = -2362.7778
$ws.Range("H91").Value = 1395.1428
$ws.Range("I91").Value = 1115
$ws.Range("J91").Value = 1550.7778
$ws.Range("K91").Value = 1115
$ws.Range("L91").Value = 1550.7778
$ws.Range("M91").Value = 289
$ws.Range("N91").Value = -4358.7778
$ws.Range("H116").Value = 9404583
$ws.Range("I116").Value = 13582678
$ws.Range("J116").Value = 3871.25
$ws.Range("K116").Value = 13582678
$ws.Range("L116").Value = 3871.25
$ws.Range("M116").Value = -13579236
$ws.Range("N116").Value = -10755.25
$ws.Range("H125").Value = 6898.923
$ws.Range("J125").Value = 4854.8887
$ws.Range("L125").Value = 43693.99830000001
$ws.Range("N125").Value = -48613.99830000001
$ws.Range("H138").Value = 317187.66
$ws.Range("I138").Value = 1098191.6
$ws.Range("J138").Value = 4786.1
$ws.Range("K138").Value = 3294574.8
$ws.Range("L138").Value = 14358.3
$ws.Range("M138").Value = -3289434.8
$ws.Range("N138").Value = -24638.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8798.102000000001
$ws.Range("J32").Value = 16998.334
$ws.Range("L32").Value = 16998.334
$ws.Range("N32").Value = -17572.334
$ws.Range("H74").Value = 1927.1052
$ws.Range("I74").Value = 554
$ws.Range("K74").Value = 554
$ws.Range("M74").Value = 320
$ws.Range("H77").Value = 1927.1052
$ws.Range("I77").Value = 554
$ws.Range("K77").Value = 2770
$ws.Range("M77").Value = 1598
$ws.Range("H122").Value = 685186.7
$ws.Range("I122").Value = 2880.35
$ws.Range("J122").Value = 7508250
$ws.Range("K122").Value = 8641.049999999999
$ws.Range("L122").Value = 22524750
$ws.Range("M122").Value = -6191.049999999999
$ws.Range("N122").Value = -22529650
$ws.Range("H132").Value = 2470.2285
$ws.Range("I132").Value = 800.4783
$ws.Range("J132").Value = 5670.5835
$ws.Range("K132").Value = 2401.4349
$ws.Range("L132").Value = 17011.7505
$ws.Range("M132").Value = 128.5650999999998
$ws.Range("N132").Value = -22071.7505

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 150
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H80").Value = 282.66666
$ws.Range("J80").Value = 231.53334
$ws.Range("L80").Value = 231.53334
$ws.Range("N80").Value = -2227.53334
$ws.Range("H83").Value = 282.66666
$ws.Range("J83").Value = 231.53334
$ws.Range("L83").Value = 1157.6667
$ws.Range("N83").Value = -11141.6667
$ws.Range("H86").Value = 5319.276
$ws.Range("I86").Value = 5932.5713
$ws.Range("J86").Value = 3709.375
$ws.Range("K86").Value = 5932.5713
$ws.Range("L86").Value = 3709.375
$ws.Range("M86").Value = -4809.5713
$ws.Range("N86").Value = -5955.375
$ws.Range("H89").Value = 5319.276
$ws.Range("I89").Value = 5932.5713
$ws.Range("J89").Value = 3709.375
$ws.Range("K89").Value = 29662.8565
$ws.Range("L89").Value = 18546.875
$ws.Range("M89").Value = -24046.8565
$ws.Range("N89").Value = -29778.875
$ws.Range("H107").Value = 7377
$ws.Range("I107").Value = 8471.25
$ws.Range("K107").Value = 8471.25
$ws.Range("M107").Value = -6551.25
$ws.Range("H134").Value = 6379.5
$ws.Range("I134").Value = 6927.7896
$ws.Range("K134").Value = 20783.3688
$ws.Range("M134").Value = -18248.3688

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 54164.633
$ws.Range("I16").Value = 1330.7858
$ws.Range("J16").Value = 202099.4
$ws.Range("K16").Value = 1330.7858
$ws.Range("L16").Value = 202099.4
$ws.Range("M16").Value = -1043.7858
$ws.Range("N16").Value = -202673.4
$ws.Range("H58").Value = 2840.9565
$ws.Range("I58").Value = 2072.5833
$ws.Range("K58").Value = 2072.5833
$ws.Range("M58").Value = -1869.5833
$ws.Range("H107").Value = 71437670
$ws.Range("I107").Value = 111124456
$ws.Range("K107").Value = 111124456
$ws.Range("M107").Value = -111122536
$ws.Range("H113").Value = 54164.633
$ws.Range("I113").Value = 1330.7858
$ws.Range("J113").Value = 202099.4
$ws.Range("K113").Value = 1330.7858
$ws.Range("L113").Value = 202099.4
$ws.Range("M113").Value = 839.2141999999999
$ws.Range("N113").Value = -206439.4
$ws.Range("H136").Value = 2840.9565
$ws.Range("I136").Value = 2072.5833
$ws.Range("K136").Value = 6217.749899999999
$ws.Range("M136").Value = -3667.749899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 228632.84
$ws.Range("I5").Value = 444.44446
$ws.Range("J5").Value = 386609.44
$ws.Range("K5").Value = 1333.33338
$ws.Range("L5").Value = 1159828.32
$ws.Range("M5").Value = -1221.33338
$ws.Range("N5").Value = -1160052.32
$ws.Range("H113").Value = 832.51514
$ws.Range("J113").Value = 877.75
$ws.Range("L113").Value = 2633.25
$ws.Range("N113").Value = -6973.25
$ws.Range("H122").Value = 4201.864
$ws.Range("J122").Value = 5343.0625
$ws.Range("L122").Value = 48087.5625
$ws.Range("N122").Value = -52987.5625
$ws.Range("H129").Value = 55557344
$ws.Range("I129").Value = 807.5
$ws.Range("J129").Value = 166670420
$ws.Range("K129").Value = 2422.5
$ws.Range("L129").Value = 500011260
$ws.Range("M129").Value = 2577.5
$ws.Range("N129").Value = -500021260
$ws.Range("H132").Value = 32102.25
$ws.Range("J132").Value = 57655.91
$ws.Range("L132").Value = 518903.1900000001
$ws.Range("N132").Value = -523963.1900000001
$ws.Range("H135").Value = 228632.84
$ws.Range("I135").Value = 444.44446
$ws.Range("J135").Value = 386609.44
$ws.Range("K135").Value = 4000.00014
$ws.Range("L135").Value = 3479484.96
$ws.Range("M135").Value = -1465.00014
$ws.Range("N135").Value = -3484554.96
$ws.Range("H137").Value = 4901.206
$ws.Range("I137").Value = 2145.4348
$ws.Range("J137").Value = 10663.272
$ws.Range("K137").Value = 6436.3044
$ws.Range("L137").Value = 31989.816
$ws.Range("M137").Value = -1336.3044
$ws.Range("N137").Value = -42189.81600000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 21236.5
$ws.Range("I80").Value = 21236.5
$ws.Range("K80").Value = 21236.5
$ws.Range("M80").Value = -20238.5
$ws.Range("H83").Value = 21236.5
$ws.Range("I83").Value = 21236.5
$ws.Range("K83").Value = 106182.5
$ws.Range("M83").Value = -101190.5
$ws.Range("H107").Value = 404.3913
$ws.Range("I107").Value = 417.2353
$ws.Range("K107").Value = 417.2353
$ws.Range("M107").Value = 1502.7647
$ws.Range("H126").Value = 8521.909
$ws.Range("I126").Value = 16647.25
$ws.Range("J126").Value = 3878.8572
$ws.Range("K126").Value = 49941.75
$ws.Range("L126").Value = 11636.5716
$ws.Range("M126").Value = -47471.75
$ws.Range("N126").Value = -16576.5716
$ws.Range("H132").Value = 2387.375
$ws.Range("I132").Value = 2589.35
$ws.Range("J132").Value = 1377.5
$ws.Range("K132").Value = 7768.049999999999
$ws.Range("L132").Value = 4132.5
$ws.Range("M132").Value = -5238.049999999999
$ws.Range("N132").Value = -9192.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 8379.4
$ws.Range("I68").Value = 2272
$ws.Range("J68").Value = 12451
$ws.Range("K68").Value = 2272
$ws.Range("L68").Value = 12451
$ws.Range("M68").Value = -1523
$ws.Range("N68").Value = -13949
$ws.Range("H71").Value = 8379.4
$ws.Range("I71").Value = 2272
$ws.Range("J71").Value = 12451
$ws.Range("K71").Value = 11360
$ws.Range("L71").Value = 62255
$ws.Range("M71").Value = -7616
$ws.Range("N71").Value = -69743
$ws.Range("H82").Value = 2150.6667
$ws.Range("I82").Value = 2159.5
$ws.Range("K82").Value = 2159.5
$ws.Range("M82").Value = -1798.5
$ws.Range("H85").Value = 2150.6667
$ws.Range("I85").Value = 2159.5
$ws.Range("K85").Value = 2159.5
$ws.Range("M85").Value = -911.5
$ws.Range("H132").Value = 440209.34
$ws.Range("I132").Value = 597396.1
$ws.Range("J132").Value = 3579.4443
$ws.Range("K132").Value = 1792188.3
$ws.Range("L132").Value = 10738.3329
$ws.Range("M132").Value = -1789658.3
$ws.Range("N132").Value = -15798.3329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 21322.459
$ws.Range("I126").Value = 28443.53
$ws.Range("J126").Value = 4028.4285
$ws.Range("K126").Value = 85330.59
$ws.Range("L126").Value = 12085.2855
$ws.Range("M126").Value = -82860.59
$ws.Range("N126").Value = -17025.2855
$ws.Range("H132").Value = 6792.6133
$ws.Range("I132").Value = 7611.4727
$ws.Range("K132").Value = 22834.4181
$ws.Range("M132").Value = -20304.4181

